# Apply crypto price/volume updates per the commit diff.
# Values that look like plain numbers (e.g. "8.37") must be forced to
# text so Excel does not silently coerce them to floating point and
# lose the original formatting (trailing zeros, etc.) - the workbook
# stores every Price/Volume cell as text (inlineStr).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.349.98"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "3.521.79"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("D8").Value = "3.517.80"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.581"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "4.101.12"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "615.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -10.30%  "
$ws.Range("D18").Value = "3.531.28"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "69.496.67"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.883"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.58%  "
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "571.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("B36").Value = "Cosmos"
$ws.Range("C36").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("E38").Value = "  -3.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0445"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "3.386.13"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("E44").Value = "  -2.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "33.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").Value = "0.0₃0705"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.95%  "
